# Insert a new weekly price record at the top of the dated listing (row 9),
# pushing the existing rows 9-39 down to rows 10-40.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(9).Insert()

$ws.Cells.Item(9, 1).Value = 9
$ws.Cells.Item(9, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(9, 3).Value = "Metropolitana"
$ws.Cells.Item(9, 4).Value = 44677
$ws.Cells.Item(9, 5).Value = 13
$ws.Cells.Item(9, 6).Value = "Fruta"
$ws.Cells.Item(9, 7).Value = 100104
$ws.Cells.Item(9, 8).Value = "Frutos de pepita"
$ws.Cells.Item(9, 9).Value = 100104003
$ws.Cells.Item(9, 10).Value = "Membrillo"
$ws.Cells.Item(9, 11).Value = "Champion"
$ws.Cells.Item(9, 12).Value = "Primera"
$ws.Cells.Item(9, 13).Value = 320
$ws.Cells.Item(9, 14).Value = 10000
$ws.Cells.Item(9, 15).Value = 10000
$ws.Cells.Item(9, 16).Value = 10000
$ws.Cells.Item(9, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(9, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(9, 19).Value = 556
$ws.Cells.Item(9, 20).Value = 18
